$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# New elevation field names for column A (rows 50-63), in the order the
# author apparently typed them (column by column) so that shared strings
# are appended in the right order.
$fieldNames = @(
    "GainCorrectedElevation",
    "GainElevation",
    "GainUncorrectedElevation",
    "LossCorrectedElevation",
    "LossElevation",
    "LossUncorrectedElevation",
    "MaxCorrectedElevation",
    "MaxElevation",
    "MaxUncorrectedElevation",
    "MinCorrectedElevation",
    "MinElevation",
    "MinUncorrectedElevation",
    "SumSampleCountElevation",
    "WeightedMeanElevation"
)

$startRow = 50
$endRow = 63

# Fill column A first (field names)
for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $r = $startRow + $i
    $ws2.Cells.Item($r, 1).Value = $fieldNames[$i]
}

# Fill column B (activity_type) - always "all"
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws2.Cells.Item($r, 2).Value = "all"
}

# Fill column C (metric) - "meterelevation" so it isn't compounded with
# the regular "meter" unit conversion used elsewhere.
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws2.Cells.Item($r, 3).Value = "meterelevation"
}

# Fill column D (statute) - "footelevation"
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws2.Cells.Item($r, 4).Value = "footelevation"
}

# Update view/selection state: gc_fields_uom (sheet2) becomes the active tab.
$ws1.Activate()
$ws1.Range("A23").Select() | Out-Null

$ws2.Activate()
$ws2.Range("C52").Select() | Out-Null
